$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.350.59"
$ws.Range("E2").Value = "  -0.08%  "
# Row 3
$ws.Range("D3").Value = "3.431.79"
$ws.Range("E3").Value = "  +1.63%  "
# Row 4
$ws.Range("E4").Value = "  -0.03%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.75"
$ws.Range("E5").Value = "  +0.29%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.81"
$ws.Range("E6").Value = "  +2.01%  "
# Row 7
$ws.Range("E7").Value = "  +0.00%  "
# Row 8
$ws.Range("D8").Value = "3.431.94"
$ws.Range("E8").Value = "  +1.65%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.472"
$ws.Range("E9").Value = "  -0.51%  "
# Row 10
$ws.Range("E10").Value = "  +3.56%  "
# Row 11
$ws.Range("E11").Value = "  -1.20%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.382"
$ws.Range("E12").Value = "  -2.54%  "
# Row 13
$ws.Range("D13").Value = "4.018.93"
$ws.Range("E13").Value = "  +1.63%  "
# Row 14
$ws.Range("E14").Value = "  -0.75%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.81"
$ws.Range("E15").Value = "  +2.95%  "
# Row 16
$ws.Range("E16").Value = "  -1.26%  "
# Row 17
$ws.Range("D17").Value = "3.427.68"
$ws.Range("E17").Value = "  +1.40%  "
# Row 18
$ws.Range("D18").Value = "61.466.30"
$ws.Range("E18").Value = "  -0.11%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.98"
$ws.Range("E19").Value = "  +1.80%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.96"
$ws.Range("E20").Value = "  -0.25%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.44"
$ws.Range("E21").Value = "  +1.11%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "384.50"
$ws.Range("E22").Value = "  +1.74%  "
# Row 23
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.556"
$ws.Range("E23").Value = "  +0.27%  "
# Row 24
$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D24").Value = "3.543.99"
$ws.Range("E24").Value = "  +0.82%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.01"
$ws.Range("E25").Value = "  +0.65%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.77"
$ws.Range("E26").Value = "  +0.03%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000123"
$ws.Range("E27").Value = "  -1.99%  "
# Row 28
$ws.Range("E28").Value = "  +8.86%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.59"
$ws.Range("E29").Value = "  -8.22%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.54"
$ws.Range("E30").Value = "  +0.68%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.35%  "
# Row 32
$ws.Range("E32").Value = "  -1.51%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.16"
$ws.Range("E33").Value = "  +0.12%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.76"
$ws.Range("E35").Value = "  +0.94%  "
# Row 36
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.16"
$ws.Range("E36").Value = "  -0.62%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.94"
$ws.Range("E37").Value = "  +1.87%  "
# Row 38
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.56"
$ws.Range("E38").Value = "  +2.16%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "165.91"
$ws.Range("E39").Value = "  +0.37%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0777"
$ws.Range("E40").Value = "  +0.76%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.40"
$ws.Range("E41").Value = "  +8.57%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.74"
$ws.Range("E42").Value = "  +0.60%  "
# Row 43
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.784"
$ws.Range("E43").Value = "  +1.32%  "
# Row 44
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.00%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.16"
$ws.Range("E45").Value = "  +1.29%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.42"
$ws.Range("E46").Value = "  +0.84%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.18"
$ws.Range("E47").Value = "  -2.51%  "
# Row 48
$ws.Range("D48").Value = "2.579.38"
$ws.Range("E48").Value = "  +9.14%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.89"
$ws.Range("E49").Value = "  +5.81%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.81"
$ws.Range("E50").Value = "  -0.31%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0262"
$ws.Range("E51").Value = "  -0.89%  "
